$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I3").Value = 1.33
$ws.Range("J3").Value = 8.5
$ws.Range("L3").Value = 1.83
$ws.Range("M3").Value = 1.04
$ws.Range("N3").Value = 13
$ws.Range("U3").Value = 2.38
$ws.Range("V3").Value = 1.53
